$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037576114560824
$ws.Range("D2").Value = 1.048704673305064
$ws.Range("E2").Value = 1.055434095325958
$ws.Range("F2").Value = 1.061749946663601
$ws.Range("I2").Value = 1.044820647590991
$ws.Range("J2").Value = 1.042678235572573
$ws.Range("K2").Value = 1.051463775017783
$ws.Range("L2").Value = 1.058174564293716
$ws.Range("M2").Value = 1.06447316064288
$ws.Range("N2").Value = 1.018173480693993
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038351130261965
$ws.Range("D3").Value = 1.049338052024739
$ws.Range("E3").Value = 1.05626950743407
$ws.Range("F3").Value = 1.0625435093216
$ws.Range("I3").Value = 1.045028726524142
$ws.Range("J3").Value = 1.043098564207513
$ws.Range("K3").Value = 1.051909715028511
$ws.Range("L3").Value = 1.05882336086082
$ws.Range("M3").Value = 1.065081458333968
$ws.Range("N3").Value = 1.018313202415649
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038853159574755
$ws.Range("D4").Value = 1.049748352348916
$ws.Range("E4").Value = 1.056811413666079
$ws.Range("F4").Value = 1.063057974178271
$ws.Range("I4").Value = 1.045162342050538
$ws.Range("J4").Value = 1.043370391676055
$ws.Range("K4").Value = 1.052198038946603
$ws.Range("L4").Value = 1.059243863185372
$ws.Range("M4").Value = 1.065475384308004
$ws.Range("N4").Value = 1.018403543808349
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039064340720739
$ws.Range("D5").Value = 1.049920951302893
$ws.Range("E5").Value = 1.05703954948522
$ws.Range("F5").Value = 1.063274486970401
$ws.Range("I5").Value = 1.04521826764603
$ws.Range("J5").Value = 1.043484629981328
$ws.Range("K5").Value = 1.052319194070575
$ws.Range("L5").Value = 1.059420805256724
$ws.Range("M5").Value = 1.065641064897633
$ws.Range("N5").Value = 1.01844150662086
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039099806419057
$ws.Range("D6").Value = 1.049949937745772
$ws.Range("E6").Value = 1.057077873106715
$ws.Range("F6").Value = 1.063310853951581
$ws.Range("I6").Value = 1.045227643326812
$ws.Range("J6").Value = 1.043503808832366
$ws.Range("K6").Value = 1.052339533209478
$ws.Range("L6").Value = 1.059450524111562
$ws.Range("M6").Value = 1.065668887693887
$ws.Range("N6").Value = 1.018447879751135
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03885598088679
$ws.Range("D7").Value = 1.049750658199859
$ws.Range("E7").Value = 1.056814460779396
$ws.Range("F7").Value = 1.063060866324715
$ws.Range("I7").Value = 1.045163090300246
$ws.Range("J7").Value = 1.043371918284641
$ws.Range("K7").Value = 1.052199658050094
$ws.Range("L7").Value = 1.059246226856074
$ws.Range("M7").Value = 1.065477597851092
$ws.Range("N7").Value = 1.018404051135531
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037837921362708
$ws.Range("D8").Value = 1.048918629824617
$ws.Range("E8").Value = 1.055716148328684
$ws.Range("F8").Value = 1.062017931391744
$ws.Range("I8").Value = 1.044891180641806
$ws.Range("J8").Value = 1.042820318573448
$ws.Range("K8").Value = 1.0516145291207
$ws.Range("L8").Value = 1.058393684739211
$ws.Range("M8").Value = 1.064678671096472
$ws.Range("N8").Value = 1.018220714092019
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036048204195047
$ws.Range("D9").Value = 1.04745611101949
$ws.Range("E9").Value = 1.053791125581502
$ws.Range("F9").Value = 1.06018771911094
$ws.Range("I9").Value = 1.044404226049586
$ws.Range("J9").Value = 1.041847211757882
$ws.Range("K9").Value = 1.050581764858183
$ws.Range("L9").Value = 1.056896738137864
$ws.Range("M9").Value = 1.063273361045808
$ws.Range("N9").Value = 1.017897150733375
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034858014114535
$ws.Range("D10").Value = 1.046483646307254
$ws.Range("E10").Value = 1.052514854899284
$ws.Range("F10").Value = 1.058972792517313
$ws.Range("I10").Value = 1.044074387550461
$ws.Range("J10").Value = 1.041197798257598
$ws.Range("K10").Value = 1.049892201516487
$ws.Range("L10").Value = 1.055902467293897
$ws.Range("M10").Value = 1.062338271155341
$ws.Range("N10").Value = 1.017681132611952
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034343372302617
$ws.Range("D11").Value = 1.046063186564389
$ws.Range("E11").Value = 1.051963920429688
$ws.Range("F11").Value = 1.058447978995299
$ws.Range("I11").Value = 1.04393034099509
$ws.Range("J11").Value = 1.040916450808582
$ws.Range("K11").Value = 1.049593380771425
$ws.Range("L11").Value = 1.055472834745045
$ws.Range("M11").Value = 1.061933811498742
$ws.Range("N11").Value = 1.017587526743305
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03415232118313
$ws.Range("D12").Value = 1.045907104768343
$ws.Range("E12").Value = 1.051759536241215
$ws.Range("F12").Value = 1.058253230972667
$ws.Range("I12").Value = 1.043876652626618
$ws.Range("J12").Value = 1.040811925025848
$ws.Range("K12").Value = 1.049482351525294
$ws.Range("L12").Value = 1.055313386051392
$ws.Range("M12").Value = 1.061783645091022
$ws.Range("N12").Value = 1.017552747509799
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034193297300189
$ws.Range("D13").Value = 1.045940580482882
$ws.Range("E13").Value = 1.05180336566127
$ws.Range("F13").Value = 1.058294996403347
$ws.Range("I13").Value = 1.04388817723389
$ws.Range("J13").Value = 1.040834347082205
$ws.Range("K13").Value = 1.049506169187326
$ws.Range("L13").Value = 1.055347582156547
$ws.Range("M13").Value = 1.061815853193792
$ws.Range("N13").Value = 1.017560208211627
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034327577707684
$ws.Range("D14").Value = 1.046050282831521
$ws.Range("E14").Value = 1.051947020705666
$ws.Range("F14").Value = 1.058431877152572
$ws.Range("I14").Value = 1.043925906827194
$ws.Range("J14").Value = 1.040907811091901
$ws.Range("K14").Value = 1.049584203742922
$ws.Range("L14").Value = 1.055459651882491
$ws.Range("M14").Value = 1.061921397301044
$ws.Range("N14").Value = 1.017584652079214
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034410326898405
$ws.Range("D15").Value = 1.04611788680479
$ws.Range("E15").Value = 1.05203556549178
$ws.Range("F15").Value = 1.05851623928434
$ws.Range("I15").Value = 1.043949129043266
$ws.Range("J15").Value = 1.040953071973509
$ws.Range("K15").Value = 1.049632278956591
$ws.Range("L15").Value = 1.055528719808584
$ws.Range("M15").Value = 1.061986435560852
$ws.Range("N15").Value = 1.017599711461824
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034892184488468
$ws.Range("D16").Value = 1.046511564146075
$ws.Range("E16").Value = 1.052551454586254
$ws.Range("F16").Value = 1.059007649342698
$ws.Range("I16").Value = 1.044083921725767
$ws.Range("J16").Value = 1.041216467371683
$ws.Range("K16").Value = 1.049912028430446
$ws.Range("L16").Value = 1.055930999561826
$ws.Range("M16").Value = 1.0623651232393
$ws.Range("N16").Value = 1.017687343519177
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035194634644844
$ws.Range("D17").Value = 1.046758675946451
$ws.Range("E17").Value = 1.052875514598545
$ws.Range("F17").Value = 1.059316236118266
$ws.Range("I17").Value = 1.044168146448984
$ws.Range("J17").Value = 1.041381649869055
$ws.Range("K17").Value = 1.050087445906304
$ws.Range("L17").Value = 1.056183579422652
$ws.Range("M17").Value = 1.062602783165835
$ws.Range("N17").Value = 1.017742294753022
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035371117890589
$ws.Range("D18").Value = 1.046902872145657
$ws.Range("E18").Value = 1.053064697122714
$ws.Range("F18").Value = 1.059496350771619
$ws.Range("I18").Value = 1.044217155118325
$ws.Range("J18").Value = 1.041477983743905
$ws.Range("K18").Value = 1.050189741156371
$ws.Range("L18").Value = 1.056330991053608
$ws.Range("M18").Value = 1.062741448558741
$ws.Range("N18").Value = 1.017774340228642
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035431305795924
$ws.Range("D19").Value = 1.046952049452527
$ws.Range("E19").Value = 1.053129231193493
$ws.Range("F19").Value = 1.059557785707996
$ws.Range("I19").Value = 1.044233845734967
$ws.Range("J19").Value = 1.041510828644155
$ws.Range("K19").Value = 1.050224617283026
$ws.Range("L19").Value = 1.056381269156025
$ws.Range("M19").Value = 1.062788737009291
$ws.Range("N19").Value = 1.017785265765854
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035162177435857
$ws.Range("D20").Value = 1.046732156961921
$ws.Range("E20").Value = 1.052840729068004
$ws.Range("F20").Value = 1.059283115125206
$ws.Range("I20").Value = 1.044159122151581
$ws.Range("J20").Value = 1.041363928818226
$ws.Range("K20").Value = 1.050068627618838
$ws.Range("L20").Value = 1.056156471078525
$ws.Range("M20").Value = 1.062577280108889
$ws.Range("N20").Value = 1.017736399690143
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03428803242972
$ws.Range("D21").Value = 1.046017975567483
$ws.Range("E21").Value = 1.051904710752934
$ws.Range("F21").Value = 1.058391563891558
$ws.Range("I21").Value = 1.043914801448788
$ws.Range("J21").Value = 1.040886178336279
$ws.Range("K21").Value = 1.049561225436038
$ws.Range("L21").Value = 1.055426646374067
$ws.Range("M21").Value = 1.061890315303201
$ws.Range("N21").Value = 1.017577454237567
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033739058741637
$ws.Range("D22").Value = 1.045569495874184
$ws.Range("E22").Value = 1.051317689198335
$ws.Range("F22").Value = 1.057832116867816
$ws.Range("I22").Value = 1.043760128607059
$ws.Range("J22").Value = 1.040585677865792
$ws.Range("K22").Value = 1.049242006022301
$ws.Range("L22").Value = 1.054968564133439
$ws.Range("M22").Value = 1.06145878707701
$ws.Range("N22").Value = 1.017477462183795
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034030019182927
$ws.Range("D23").Value = 1.045807190238291
$ws.Range("E23").Value = 1.051628738457823
$ws.Range("F23").Value = 1.058128584714304
$ws.Range("I23").Value = 1.043842223668954
$ws.Range("J23").Value = 1.040744989774186
$ws.Range("K23").Value = 1.049411248329375
$ws.Range("L23").Value = 1.055211326961962
$ws.Range("M23").Value = 1.061687510442544
$ws.Range("N23").Value = 1.017530475088243
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035176843237952
$ws.Range("D24").Value = 1.046744139564193
$ws.Range("E24").Value = 1.052856446646091
$ws.Range("F24").Value = 1.05929808070206
$ws.Range("I24").Value = 1.044163200208417
$ws.Range("J24").Value = 1.041371936242593
$ws.Range("K24").Value = 1.050077130863539
$ws.Range("L24").Value = 1.056168719906654
$ws.Range("M24").Value = 1.062588803711841
$ws.Range("N24").Value = 1.017739063435912
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.036510376118816
$ws.Range("D25").Value = 1.047833766373925
$ws.Range("E25").Value = 1.054287551595205
$ws.Range("F25").Value = 1.060659962743335
$ws.Range("I25").Value = 1.0445310362231
$ws.Range("J25").Value = 1.042098907827693
$ws.Range("K25").Value = 1.050848951067897
$ws.Range("L25").Value = 1.057283091183278
$ws.Range("M25").Value = 1.063636360492799
$ws.Range("N25").Value = 1.017980856093774
